$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update Marking -> Right count (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update Total -> Right count (B12): 48 -> 80
$ws.Range("B12").Value = 80

# Update Total -> Max (E12): "47/84" -> "80/140"
$ws.Range("E12").Value = "80/140"
